$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F ("想去人数")
$updates = @{
    2  = 8464
    3  = 8146
    4  = 139
    9  = 144
    10 = 198
    12 = 734
    14 = 4225
    16 = 76
    17 = 18
    19 = 149
    20 = 116
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
